$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The latest GSC export dropped the oldest date in the rolling window
# (2025-10-10, which had no coverage data yet). Remove that leading row
# so the remaining dates/values shift up to match the refreshed export.
$ws.Rows.Item(2).Delete()
